$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# Replace LOOKUP formulas (referencing the external workbook) with their static cached values.
# This removes the dependency on the external link so it can be dropped from the workbook.
$ws.Range("B2").Value = 'Realizar el lanzamiento del ciclo #2 de TSPi.'
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("B3").Value = 'Definir la estrategía de desarrolo del ciclo #2 de TSPi.'
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 4
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("B4").Value = 'Elaborar el plan del ciclo #2 de TSPi.'
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0
$ws.Range("B5").Value = 'Crear el esquema del documento de arquitectura.'
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 4
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("B6").Value = 'Reunión de equipo para analizar la versión final del documento de requerimientos.'
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 7.5
$ws.Range("F6").Value = 4
$ws.Range("H6").Value = 1.5
$ws.Range("I6").Value = 1.5
$ws.Range("J6").Value = 1.5
$ws.Range("K6").Value = 1.5
$ws.Range("L6").Value = 1.5
$ws.Range("B7").Value = 'Crear la agenda para la reunión #3 con el cliente.'
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 4
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.5
$ws.Range("L7").Value = 0
$ws.Range("B8").Value = 'Reunión #3 con el cliente.'
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0
$ws.Range("B9").Value = 'Prepararse para presentar al equipo la herramienta RedMine.'
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("B10").Value = 'Presentar al equipo de la herramienta Redmine.'
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 4
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("B11").Value = 'Elaborar la introducción del documento de arquitectura.'
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 4
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("B12").Value = 'Elaborar el fondo del documento de arquitectura.'
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 4
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("B13").Value = 'Elaborar el diagrama de contexto de la arquitectura.'
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 4
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("B14").Value = 'Documentar las tácticas y patrones a utilizar en el documento de arquitectura.'
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("B15").Value = 'Elaborar el modelo físico de data de la arquitectura.'
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("B16").Value = 'Elaborar el diagrama de flujo del algoritmo de calendarización.'
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2
$ws.Range("B17").Value = 'Elaborar el diagrama de la estructura de archivos del plug-in.'
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 5
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2
$ws.Range("B18").Value = 'Elaborar la conclusión del documento de arquitectura.'
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.5
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("B19").Value = 'Cursar el tutorial básico de ruby.'
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 5
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0
$ws.Range("B20").Value = 'Prepararse para presentar al equipo el framework Rails.'
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3
$ws.Range("B21").Value = 'Presentar al equipo el framework Rails.'
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 5
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("B22").Value = 'Mockup de la vista de la calendarización.'
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 2.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2.5
$ws.Range("L22").Value = 0
$ws.Range("B23").Value = 'Elaborar el reporte de cierre del ciclo #2 de TSPi.'
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 1
